$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": row 3 corresponds to the fa9f417a file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-07-17 11:07:44"

# --- Sheet "zh-cn": row 3 corresponds to the fa9f417a file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-17 11:07:41"

# --- Sheet "de-de": row 3 corresponds to the fa9f417a file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-17 11:07:44"
